$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$d.Paragraphs.Item(1).Range.Text = "2023-10-10 Tuesday"

# Update each answer cell in the 20x5 table, addressed by (row, col)
# to avoid ambiguity from duplicate cell text appearing more than once.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "72-26=46"
$t.Cell(1, 2).Range.Text = "18+23=41"
$t.Cell(1, 3).Range.Text = "50+16=66"
$t.Cell(1, 4).Range.Text = "90-72=18"
$t.Cell(1, 5).Range.Text = "85-15=70"
$t.Cell(2, 1).Range.Text = "52+21=73"
$t.Cell(2, 2).Range.Text = "4+90=94"
$t.Cell(2, 3).Range.Text = "82-2=80"
$t.Cell(2, 4).Range.Text = "98-66=32"
$t.Cell(2, 5).Range.Text = "51+16=67"
$t.Cell(3, 1).Range.Text = "51+24=75"
$t.Cell(3, 2).Range.Text = "86+13=99"
$t.Cell(3, 3).Range.Text = "81-32=49"
$t.Cell(3, 4).Range.Text = "7+23=30"
$t.Cell(3, 5).Range.Text = "30+8=38"
$t.Cell(4, 1).Range.Text = "98-70=28"
$t.Cell(4, 2).Range.Text = "91-70=21"
$t.Cell(4, 3).Range.Text = "30+68=98"
$t.Cell(4, 4).Range.Text = "7+16=23"
$t.Cell(4, 5).Range.Text = "63-22=41"
$t.Cell(5, 1).Range.Text = "97-84=13"
$t.Cell(5, 2).Range.Text = "64+30=94"
$t.Cell(5, 3).Range.Text = "39+46=85"
$t.Cell(5, 4).Range.Text = "84-46=38"
$t.Cell(5, 5).Range.Text = "64-24=40"
$t.Cell(6, 1).Range.Text = "64-57=7"
$t.Cell(6, 2).Range.Text = "20+51=71"
$t.Cell(6, 3).Range.Text = "92-84=8"
$t.Cell(6, 4).Range.Text = "3+90=93"
$t.Cell(6, 5).Range.Text = "5+86=91"
$t.Cell(7, 1).Range.Text = "76-15=61"
$t.Cell(7, 2).Range.Text = "15+40=55"
$t.Cell(7, 3).Range.Text = "66+1=67"
$t.Cell(7, 4).Range.Text = "1+16=17"
$t.Cell(7, 5).Range.Text = "66-49=17"
$t.Cell(8, 1).Range.Text = "2+86=88"
$t.Cell(8, 2).Range.Text = "30+63=93"
$t.Cell(8, 3).Range.Text = "58-26=32"
$t.Cell(8, 4).Range.Text = "4+72=76"
$t.Cell(8, 5).Range.Text = "57+20=77"
$t.Cell(9, 1).Range.Text = "2+31=33"
$t.Cell(9, 2).Range.Text = "75-28=47"
$t.Cell(9, 3).Range.Text = "49-8=41"
$t.Cell(9, 4).Range.Text = "74-34=40"
$t.Cell(9, 5).Range.Text = "9+50=59"
$t.Cell(10, 1).Range.Text = "6+58=64"
$t.Cell(10, 2).Range.Text = "65-14=51"
$t.Cell(10, 3).Range.Text = "50-15=35"
$t.Cell(10, 4).Range.Text = "21+5=26"
$t.Cell(10, 5).Range.Text = "27-18=9"
$t.Cell(11, 1).Range.Text = "27+17=44"
$t.Cell(11, 2).Range.Text = "46+14=60"
$t.Cell(11, 3).Range.Text = "97-37=60"
$t.Cell(11, 4).Range.Text = "93-13=80"
$t.Cell(11, 5).Range.Text = "15-13=2"
$t.Cell(12, 1).Range.Text = "39+33=72"
$t.Cell(12, 2).Range.Text = "11+53=64"
$t.Cell(12, 3).Range.Text = "18+44=62"
$t.Cell(12, 4).Range.Text = "1+18=19"
$t.Cell(12, 5).Range.Text = "34-3=31"
$t.Cell(13, 1).Range.Text = "60+8=68"
$t.Cell(13, 2).Range.Text = "21+75=96"
$t.Cell(13, 3).Range.Text = "29+67=96"
$t.Cell(13, 4).Range.Text = "81-38=43"
$t.Cell(13, 5).Range.Text = "58-56=2"
$t.Cell(14, 1).Range.Text = "66-13=53"
$t.Cell(14, 2).Range.Text = "14+27=41"
$t.Cell(14, 3).Range.Text = "14+39=53"
$t.Cell(14, 4).Range.Text = "52-40=12"
$t.Cell(14, 5).Range.Text = "45+40=85"
$t.Cell(15, 1).Range.Text = "19+43=62"
$t.Cell(15, 2).Range.Text = "31+64=95"
$t.Cell(15, 3).Range.Text = "57+19=76"
$t.Cell(15, 4).Range.Text = "48-27=21"
$t.Cell(15, 5).Range.Text = "74-72=2"
$t.Cell(16, 1).Range.Text = "68+29=97"
$t.Cell(16, 2).Range.Text = "7+0=7"
$t.Cell(16, 3).Range.Text = "0+17=17"
$t.Cell(16, 4).Range.Text = "59+23=82"
$t.Cell(16, 5).Range.Text = "46-36=10"
$t.Cell(17, 1).Range.Text = "29+57=86"
$t.Cell(17, 2).Range.Text = "27+52=79"
$t.Cell(17, 3).Range.Text = "67-38=29"
$t.Cell(17, 4).Range.Text = "22+40=62"
$t.Cell(17, 5).Range.Text = "42+11=53"
$t.Cell(18, 1).Range.Text = "55+11=66"
$t.Cell(18, 2).Range.Text = "51-0=51"
$t.Cell(18, 3).Range.Text = "50+38=88"
$t.Cell(18, 4).Range.Text = "97-14=83"
$t.Cell(18, 5).Range.Text = "81-22=59"
$t.Cell(19, 1).Range.Text = "75+17=92"
$t.Cell(19, 2).Range.Text = "34+3=37"
$t.Cell(19, 3).Range.Text = "47+14=61"
$t.Cell(19, 4).Range.Text = "92-50=42"
$t.Cell(19, 5).Range.Text = "42+21=63"
$t.Cell(20, 1).Range.Text = "66-45=21"
$t.Cell(20, 2).Range.Text = "34+3=37"
$t.Cell(20, 3).Range.Text = "81-16=65"
$t.Cell(20, 4).Range.Text = "11-10=1"
$t.Cell(20, 5).Range.Text = "90-4=86"
